$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Neo4j query for the "FilesTab" row: the `File Type` and `Breed` columns
# are dropped from the RETURN clause (bento object repository revisited).
$newFilesQuery = "MATCH (f:file)-->(parent)`nWITH DISTINCT f, parent`nMATCH (f)-[*]->(c:case)<--(demo:demographic)`n MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)`n MATCH (samp:sample)-->(c) `n WHERE labels(parent)[0] IN [`"diagnosis`"]  `nWITH DISTINCT f, parent, c, demo, diag, s`nRETURN coalesce(f.file_name, '') AS ``File Name``, `n        coalesce(labels(parent)[0], '') AS ``Association``,`n        coalesce(f.file_description, '') AS ``Description``,`n        coalesce(f.file_format, '') AS ``Format``,`n        coalesce(f.file_size, '') AS ``Size``,`n        coalesce(c.case_id, '') AS ``Case ID``, `n        coalesce(diag.disease_term,'') AS Diagnosis , `n        coalesce(s.clinical_study_designation,'') AS ``Study Code``"

# Find the row whose TabName (column A) is "FilesTab" and rewrite its query
# (column B), instead of assuming a fixed row number.
$filesRow = 0
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 1; $r -le $lastRow; $r++) {
    if ($ws.Cells.Item($r, 1).Value2 -eq "FilesTab") {
        $filesRow = $r
        break
    }
}
if ($filesRow -eq 0) { $filesRow = 4 }

$ws.Cells.Item($filesRow, 2).Value = $newFilesQuery

[void]$ws.Cells.Item($filesRow, 2).Select()
$excel.ActiveWindow.ScrollRow = $filesRow
